$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark that sits after the "Purpose of this
#    document" closing paragraph (Word moves this bookmark whenever the
#    user's cursor was last positioned - here it is relocated later in the
#    document, inside the "Business Case" section).
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# 2. Insert the new "_GoBack" bookmark right after "...towards that target"
#    (before the trailing comma/period), splitting that run in two.
$rng = $d.Content
$rng.Find.Execute("target,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($rng.End - 1, $rng.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3. Trim the trailing clause, leaving a single period after the bookmark.
$old2 = ", and display highest expenditure categories, so that expenses can be reduced."
$new2 = "."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null
